$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 54: copy the formatting (number format, font, border, alignment) from the
# last existing data row (53) in column A, which carries the date style, then fill values.
$ws.Range("A54").Value = 45986
$ws.Range("A53").Copy()
$ws.Range("A54").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B54").Value = 2025
$ws.Range("C54").Value = -2.451276118722334
$ws.Range("D54").Value = 2026
$ws.Range("E54").Value = -1.596682557877005
